# Update the worksheet date and regenerate all 100 arithmetic answers
# in the table (each old expression is unique in the document, so a
# straightforward literal Find/Replace per cell is safe and unambiguous).
$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-04-01 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-02 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("92-53=39", $true, $false, $false, $false, $false, $true, 1, $false, "12+16=28", 2) | Out-Null
$d.Content.Find.Execute("35-23=12", $true, $false, $false, $false, $false, $true, 1, $false, "20-17=3", 2) | Out-Null
$d.Content.Find.Execute("0+53=53", $true, $false, $false, $false, $false, $true, 1, $false, "72-48=24", 2) | Out-Null
$d.Content.Find.Execute("75-60=15", $true, $false, $false, $false, $false, $true, 1, $false, "53-30=23", 2) | Out-Null
$d.Content.Find.Execute("75-68=7", $true, $false, $false, $false, $false, $true, 1, $false, "11+62=73", 2) | Out-Null
$d.Content.Find.Execute("21+53=74", $true, $false, $false, $false, $false, $true, 1, $false, "35+32=67", 2) | Out-Null
$d.Content.Find.Execute("91-45=46", $true, $false, $false, $false, $false, $true, 1, $false, "40+50=90", 2) | Out-Null
$d.Content.Find.Execute("11+70=81", $true, $false, $false, $false, $false, $true, 1, $false, "23+61=84", 2) | Out-Null
$d.Content.Find.Execute("82-69=13", $true, $false, $false, $false, $false, $true, 1, $false, "55-25=30", 2) | Out-Null
$d.Content.Find.Execute("67-0=67", $true, $false, $false, $false, $false, $true, 1, $false, "91-57=34", 2) | Out-Null
$d.Content.Find.Execute("54+15=69", $true, $false, $false, $false, $false, $true, 1, $false, "69+26=95", 2) | Out-Null
$d.Content.Find.Execute("46-5=41", $true, $false, $false, $false, $false, $true, 1, $false, "67+30=97", 2) | Out-Null
$d.Content.Find.Execute("51-49=2", $true, $false, $false, $false, $false, $true, 1, $false, "19+49=68", 2) | Out-Null
$d.Content.Find.Execute("60-52=8", $true, $false, $false, $false, $false, $true, 1, $false, "95-94=1", 2) | Out-Null
$d.Content.Find.Execute("68+18=86", $true, $false, $false, $false, $false, $true, 1, $false, "67+31=98", 2) | Out-Null
$d.Content.Find.Execute("90-40=50", $true, $false, $false, $false, $false, $true, 1, $false, "19-0=19", 2) | Out-Null
$d.Content.Find.Execute("69-65=4", $true, $false, $false, $false, $false, $true, 1, $false, "82-32=50", 2) | Out-Null
$d.Content.Find.Execute("25-0=25", $true, $false, $false, $false, $false, $true, 1, $false, "29+39=68", 2) | Out-Null
$d.Content.Find.Execute("93-57=36", $true, $false, $false, $false, $false, $true, 1, $false, "95-36=59", 2) | Out-Null
$d.Content.Find.Execute("62-21=41", $true, $false, $false, $false, $false, $true, 1, $false, "76-36=40", 2) | Out-Null
$d.Content.Find.Execute("54-11=43", $true, $false, $false, $false, $false, $true, 1, $false, "36+28=64", 2) | Out-Null
$d.Content.Find.Execute("57-24=33", $true, $false, $false, $false, $false, $true, 1, $false, "92-71=21", 2) | Out-Null
$d.Content.Find.Execute("53+24=77", $true, $false, $false, $false, $false, $true, 1, $false, "92-10=82", 2) | Out-Null
$d.Content.Find.Execute("31+38=69", $true, $false, $false, $false, $false, $true, 1, $false, "25+56=81", 2) | Out-Null
$d.Content.Find.Execute("5+62=67", $true, $false, $false, $false, $false, $true, 1, $false, "33+14=47", 2) | Out-Null
$d.Content.Find.Execute("68-36=32", $true, $false, $false, $false, $false, $true, 1, $false, "46+48=94", 2) | Out-Null
$d.Content.Find.Execute("75-7=68", $true, $false, $false, $false, $false, $true, 1, $false, "29-15=14", 2) | Out-Null
$d.Content.Find.Execute("46+5=51", $true, $false, $false, $false, $false, $true, 1, $false, "43-19=24", 2) | Out-Null
$d.Content.Find.Execute("8+16=24", $true, $false, $false, $false, $false, $true, 1, $false, "53-53=0", 2) | Out-Null
$d.Content.Find.Execute("69-55=14", $true, $false, $false, $false, $false, $true, 1, $false, "7+33=40", 2) | Out-Null
$d.Content.Find.Execute("49-15=34", $true, $false, $false, $false, $false, $true, 1, $false, "13-12=1", 2) | Out-Null
$d.Content.Find.Execute("82-79=3", $true, $false, $false, $false, $false, $true, 1, $false, "11+17=28", 2) | Out-Null
$d.Content.Find.Execute("54+33=87", $true, $false, $false, $false, $false, $true, 1, $false, "21+60=81", 2) | Out-Null
$d.Content.Find.Execute("8+8=16", $true, $false, $false, $false, $false, $true, 1, $false, "67-59=8", 2) | Out-Null
$d.Content.Find.Execute("11+29=40", $true, $false, $false, $false, $false, $true, 1, $false, "19+51=70", 2) | Out-Null
$d.Content.Find.Execute("97-48=49", $true, $false, $false, $false, $false, $true, 1, $false, "80-78=2", 2) | Out-Null
$d.Content.Find.Execute("37-15=22", $true, $false, $false, $false, $false, $true, 1, $false, "78-7=71", 2) | Out-Null
$d.Content.Find.Execute("70-39=31", $true, $false, $false, $false, $false, $true, 1, $false, "84+9=93", 2) | Out-Null
$d.Content.Find.Execute("9+29=38", $true, $false, $false, $false, $false, $true, 1, $false, "62-52=10", 2) | Out-Null
$d.Content.Find.Execute("30+23=53", $true, $false, $false, $false, $false, $true, 1, $false, "73-44=29", 2) | Out-Null
$d.Content.Find.Execute("42+47=89", $true, $false, $false, $false, $false, $true, 1, $false, "79-31=48", 2) | Out-Null
$d.Content.Find.Execute("94-88=6", $true, $false, $false, $false, $false, $true, 1, $false, "34-9=25", 2) | Out-Null
$d.Content.Find.Execute("65-65=0", $true, $false, $false, $false, $false, $true, 1, $false, "44-22=22", 2) | Out-Null
$d.Content.Find.Execute("16+63=79", $true, $false, $false, $false, $false, $true, 1, $false, "38+48=86", 2) | Out-Null
$d.Content.Find.Execute("1+76=77", $true, $false, $false, $false, $false, $true, 1, $false, "72-35=37", 2) | Out-Null
$d.Content.Find.Execute("59+0=59", $true, $false, $false, $false, $false, $true, 1, $false, "53+37=90", 2) | Out-Null
$d.Content.Find.Execute("32+57=89", $true, $false, $false, $false, $false, $true, 1, $false, "83-15=68", 2) | Out-Null
$d.Content.Find.Execute("58-44=14", $true, $false, $false, $false, $false, $true, 1, $false, "96-63=33", 2) | Out-Null
$d.Content.Find.Execute("49+19=68", $true, $false, $false, $false, $false, $true, 1, $false, "92-51=41", 2) | Out-Null
$d.Content.Find.Execute("53+43=96", $true, $false, $false, $false, $false, $true, 1, $false, "37+0=37", 2) | Out-Null
$d.Content.Find.Execute("1+9=10", $true, $false, $false, $false, $false, $true, 1, $false, "44+37=81", 2) | Out-Null
$d.Content.Find.Execute("14+14=28", $true, $false, $false, $false, $false, $true, 1, $false, "43-30=13", 2) | Out-Null
$d.Content.Find.Execute("24+62=86", $true, $false, $false, $false, $false, $true, 1, $false, "67+4=71", 2) | Out-Null
$d.Content.Find.Execute("78-6=72", $true, $false, $false, $false, $false, $true, 1, $false, "90-77=13", 2) | Out-Null
$d.Content.Find.Execute("87-16=71", $true, $false, $false, $false, $false, $true, 1, $false, "32+52=84", 2) | Out-Null
$d.Content.Find.Execute("51-1=50", $true, $false, $false, $false, $false, $true, 1, $false, "57-37=20", 2) | Out-Null
$d.Content.Find.Execute("82-66=16", $true, $false, $false, $false, $false, $true, 1, $false, "43+12=55", 2) | Out-Null
$d.Content.Find.Execute("84-45=39", $true, $false, $false, $false, $false, $true, 1, $false, "50-17=33", 2) | Out-Null
$d.Content.Find.Execute("62-44=18", $true, $false, $false, $false, $false, $true, 1, $false, "90-69=21", 2) | Out-Null
$d.Content.Find.Execute("30+57=87", $true, $false, $false, $false, $false, $true, 1, $false, "55-49=6", 2) | Out-Null
$d.Content.Find.Execute("81-48=33", $true, $false, $false, $false, $false, $true, 1, $false, "63-46=17", 2) | Out-Null
$d.Content.Find.Execute("72-31=41", $true, $false, $false, $false, $false, $true, 1, $false, "24+58=82", 2) | Out-Null
$d.Content.Find.Execute("78-11=67", $true, $false, $false, $false, $false, $true, 1, $false, "41-10=31", 2) | Out-Null
$d.Content.Find.Execute("95-52=43", $true, $false, $false, $false, $false, $true, 1, $false, "55-42=13", 2) | Out-Null
$d.Content.Find.Execute("52+1=53", $true, $false, $false, $false, $false, $true, 1, $false, "26+41=67", 2) | Out-Null
$d.Content.Find.Execute("40-15=25", $true, $false, $false, $false, $false, $true, 1, $false, "6+58=64", 2) | Out-Null
$d.Content.Find.Execute("28+8=36", $true, $false, $false, $false, $false, $true, 1, $false, "16-10=6", 2) | Out-Null
$d.Content.Find.Execute("80-76=4", $true, $false, $false, $false, $false, $true, 1, $false, "44-33=11", 2) | Out-Null
$d.Content.Find.Execute("1+17=18", $true, $false, $false, $false, $false, $true, 1, $false, "80+0=80", 2) | Out-Null
$d.Content.Find.Execute("11-7=4", $true, $false, $false, $false, $false, $true, 1, $false, "24+59=83", 2) | Out-Null
$d.Content.Find.Execute("94-69=25", $true, $false, $false, $false, $false, $true, 1, $false, "12+43=55", 2) | Out-Null
$d.Content.Find.Execute("3+91=94", $true, $false, $false, $false, $false, $true, 1, $false, "12+66=78", 2) | Out-Null
$d.Content.Find.Execute("60+22=82", $true, $false, $false, $false, $false, $true, 1, $false, "2+4=6", 2) | Out-Null
$d.Content.Find.Execute("67-7=60", $true, $false, $false, $false, $false, $true, 1, $false, "36-31=5", 2) | Out-Null
$d.Content.Find.Execute("33+0=33", $true, $false, $false, $false, $false, $true, 1, $false, "89-3=86", 2) | Out-Null
$d.Content.Find.Execute("4+90=94", $true, $false, $false, $false, $false, $true, 1, $false, "4+82=86", 2) | Out-Null
$d.Content.Find.Execute("42+13=55", $true, $false, $false, $false, $false, $true, 1, $false, "96-38=58", 2) | Out-Null
$d.Content.Find.Execute("26+66=92", $true, $false, $false, $false, $false, $true, 1, $false, "91-87=4", 2) | Out-Null
$d.Content.Find.Execute("54+45=99", $true, $false, $false, $false, $false, $true, 1, $false, "75-39=36", 2) | Out-Null
$d.Content.Find.Execute("1+90=91", $true, $false, $false, $false, $false, $true, 1, $false, "78+5=83", 2) | Out-Null
$d.Content.Find.Execute("13+10=23", $true, $false, $false, $false, $false, $true, 1, $false, "8-7=1", 2) | Out-Null
$d.Content.Find.Execute("57-56=1", $true, $false, $false, $false, $false, $true, 1, $false, "74+14=88", 2) | Out-Null
$d.Content.Find.Execute("24+10=34", $true, $false, $false, $false, $false, $true, 1, $false, "15+42=57", 2) | Out-Null
$d.Content.Find.Execute("15+16=31", $true, $false, $false, $false, $false, $true, 1, $false, "31+35=66", 2) | Out-Null
$d.Content.Find.Execute("26+5=31", $true, $false, $false, $false, $false, $true, 1, $false, "90-41=49", 2) | Out-Null
$d.Content.Find.Execute("65-63=2", $true, $false, $false, $false, $false, $true, 1, $false, "71-65=6", 2) | Out-Null
$d.Content.Find.Execute("36+49=85", $true, $false, $false, $false, $false, $true, 1, $false, "33+29=62", 2) | Out-Null
$d.Content.Find.Execute("57+22=79", $true, $false, $false, $false, $false, $true, 1, $false, "21+42=63", 2) | Out-Null
$d.Content.Find.Execute("29+37=66", $true, $false, $false, $false, $false, $true, 1, $false, "45+54=99", 2) | Out-Null
$d.Content.Find.Execute("10+65=75", $true, $false, $false, $false, $false, $true, 1, $false, "36+24=60", 2) | Out-Null
$d.Content.Find.Execute("15-11=4", $true, $false, $false, $false, $false, $true, 1, $false, "22+65=87", 2) | Out-Null
$d.Content.Find.Execute("66+27=93", $true, $false, $false, $false, $false, $true, 1, $false, "92-62=30", 2) | Out-Null
$d.Content.Find.Execute("41-28=13", $true, $false, $false, $false, $false, $true, 1, $false, "54-52=2", 2) | Out-Null
$d.Content.Find.Execute("67-37=30", $true, $false, $false, $false, $false, $true, 1, $false, "79-17=62", 2) | Out-Null
$d.Content.Find.Execute("0+7=7", $true, $false, $false, $false, $false, $true, 1, $false, "32+20=52", 2) | Out-Null
$d.Content.Find.Execute("73-55=18", $true, $false, $false, $false, $false, $true, 1, $false, "97-61=36", 2) | Out-Null
$d.Content.Find.Execute("55-37=18", $true, $false, $false, $false, $false, $true, 1, $false, "47+33=80", 2) | Out-Null
$d.Content.Find.Execute("49+3=52", $true, $false, $false, $false, $false, $true, 1, $false, "7+92=99", 2) | Out-Null
$d.Content.Find.Execute("50+38=88", $true, $false, $false, $false, $false, $true, 1, $false, "34+55=89", 2) | Out-Null
$d.Content.Find.Execute("36-15=21", $true, $false, $false, $false, $false, $true, 1, $false, "12+64=76", 2) | Out-Null
